$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1875
$ws.Range("C2").Value = 0.5386904761904762
$ws.Range("J2").Value = 0.008928571428571428
$ws.Range("O2").Value = 0.002976190476190476
$ws.Range("P2").Value = 0.1160714285714286
$ws.Range("S2").Value = 0.1458333333333333
$ws.Range("B3").Value = 0.02150537634408602
$ws.Range("C3").Value = 0.03225806451612903
$ws.Range("J3").Value = 0.02150537634408602
$ws.Range("P3").Value = 0.7096774193548387
$ws.Range("S3").Value = 0.2150537634408602
$ws.Range("P4").Value = 0.6458333333333334
$ws.Range("S4").Value = 0.3541666666666667
$ws.Range("B6").Value = 0.06382978723404255
$ws.Range("D6").Value = 0.02978723404255319
$ws.Range("F6").Value = 0.07234042553191489
$ws.Range("J6").Value = 0.2297872340425532
$ws.Range("O6").Value = 0.01276595744680851
$ws.Range("Q6").Value = 0.1829787234042553
$ws.Range("R6").Value = 0.06382978723404255
$ws.Range("S6").Value = 0.3446808510638298
$ws.Range("B7").Value = 0.1224489795918367
$ws.Range("D7").Value = 0.02040816326530612
$ws.Range("F7").Value = 0.04591836734693878
$ws.Range("J7").Value = 0.1479591836734694
$ws.Range("O7").Value = 0.02551020408163265
$ws.Range("Q7").Value = 0.1224489795918367
$ws.Range("R7").Value = 0.06122448979591837
$ws.Range("S7").Value = 0.4540816326530612
$ws.Range("B8").Value = 0.1064814814814815
$ws.Range("D8").Value = 0.0162037037037037
$ws.Range("F8").Value = 0.05555555555555555
$ws.Range("J8").Value = 0.1180555555555556
$ws.Range("O8").Value = 0.02546296296296296
$ws.Range("Q8").Value = 0.224537037037037
$ws.Range("R8").Value = 0.05092592592592592
$ws.Range("S8").Value = 0.4027777777777778
$ws.Range("B9").Value = 0.09195402298850575
$ws.Range("D9").Value = 0.01724137931034483
$ws.Range("F9").Value = 0.06321839080459771
$ws.Range("J9").Value = 0.1149425287356322
$ws.Range("O9").Value = 0.03448275862068965
$ws.Range("Q9").Value = 0.1781609195402299
$ws.Range("R9").Value = 0.08045977011494253
$ws.Range("S9").Value = 0.4195402298850575
$ws.Range("B10").Value = 0.122962962962963
$ws.Range("D10").Value = 0.02074074074074074
$ws.Range("E10").Value = 0.001481481481481481
$ws.Range("F10").Value = 0.06962962962962962
$ws.Range("J10").Value = 0.12
$ws.Range("O10").Value = 0.01259259259259259
$ws.Range("Q10").Value = 0.2259259259259259
$ws.Range("R10").Value = 0.05555555555555555
$ws.Range("S10").Value = 0.3711111111111111
$ws.Range("F11").Value = 0.009259259259259259
$ws.Range("G11").Value = 0.1296296296296296
$ws.Range("J11").Value = 0.09876543209876543
$ws.Range("K11").Value = 0.2006172839506173
$ws.Range("L11").Value = 0.5462962962962963
$ws.Range("S11").Value = 0.0154320987654321
$ws.Range("G12").Value = 0.7252747252747253
$ws.Range("J12").Value = 0.2142857142857143
$ws.Range("K12").Value = 0.02197802197802198
$ws.Range("L12").Value = 0.02197802197802198
$ws.Range("S12").Value = 0.01648351648351648
$ws.Range("G13").Value = 0.5490196078431373
$ws.Range("J13").Value = 0.3529411764705883
$ws.Range("S13").Value = 0.09803921568627451
$ws.Range("F15").Value = 0.04680851063829787
$ws.Range("H15").Value = 0.2042553191489362
$ws.Range("I15").Value = 0.05531914893617021
$ws.Range("J15").Value = 0.3361702127659574
$ws.Range("K15").Value = 0.05531914893617021
$ws.Range("M15").Value = 0.02127659574468085
$ws.Range("O15").Value = 0.04680851063829787
$ws.Range("S15").Value = 0.2340425531914894
$ws.Range("F16").Value = 0.005025125628140704
$ws.Range("H16").Value = 0.1658291457286432
$ws.Range("I16").Value = 0.07537688442211055
$ws.Range("J16").Value = 0.3517587939698493
$ws.Range("K16").Value = 0.1306532663316583
$ws.Range("M16").Value = 0.03517587939698492
$ws.Range("O16").Value = 0.05025125628140704
$ws.Range("S16").Value = 0.185929648241206
$ws.Range("F17").Value = 0.01807228915662651
$ws.Range("H17").Value = 0.1566265060240964
$ws.Range("I17").Value = 0.06827309236947791
$ws.Range("J17").Value = 0.463855421686747
$ws.Range("K17").Value = 0.09236947791164658
$ws.Range("M17").Value = 0.02008032128514056
$ws.Range("N17").Value = 0.002008032128514056
$ws.Range("O17").Value = 0.06425702811244979
$ws.Range("S17").Value = 0.1144578313253012
$ws.Range("F18").Value = 0.0145985401459854
$ws.Range("H18").Value = 0.1678832116788321
$ws.Range("I18").Value = 0.05109489051094891
$ws.Range("J18").Value = 0.489051094890511
$ws.Range("K18").Value = 0.08759124087591241
$ws.Range("M18").Value = 0.0145985401459854
$ws.Range("O18").Value = 0.08029197080291971
$ws.Range("S18").Value = 0.0948905109489051
$ws.Range("F19").Value = 0.0213903743315508
$ws.Range("H19").Value = 0.1902215431627196
$ws.Range("I19").Value = 0.07944996180290298
$ws.Range("J19").Value = 0.387318563789152
$ws.Range("K19").Value = 0.1184110007639419
$ws.Range("M19").Value = 0.02291825821237586
$ws.Range("N19").Value = 0.0007639419404125286
$ws.Range("O19").Value = 0.07792207792207792
$ws.Range("S19").Value = 0.1016042780748663
